# #remove static dynamic excel read
# Sheet1!C1 held a static numeric literal (567). Replace it with the
# textual value "567abc" so downstream reads pick up the dynamic/string
# content instead of the old hard-coded number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").Value = "567abc"
